$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 123
$ws.Range("A2").Value = 4563
$ws.Range("A3").Value = 89

$ws.Range("A4").Select()
